# Apply the "price refresh" update to the MY_COLLECTION_PRICES worksheet.
#
# The Date / Prices / Notes columns are stored as plain text (inline strings)
# in the workbook, not as real dates/numbers. A naive `.Value = "..."`
# assignment would let Excel's type inference turn a date-looking string
# into a date serial, or a number-looking string into a numeric value (and
# would also stamp a number-format style onto the cell). To avoid that we
# force each target cell to Text format before assigning, then reset the
# cell style back to "Normal" afterwards so no stray formatting is left
# behind on the cell (only the text value itself changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Sol Ring
Set-TextValue $ws.Range("A2") "2025-12-13"
Set-TextValue $ws.Range("E2") "2.36"
Set-TextValue $ws.Range("F2") "x2 (P/L: `$-0.64)"

# Row 3 - Dark Magician
Set-TextValue $ws.Range("A3") "2025-12-13"
Set-TextValue $ws.Range("E3") "0.25"
Set-TextValue $ws.Range("F3") " (P/L: `$-4.75)"

# Row 4 - Blue-Eyes White Dragon
Set-TextValue $ws.Range("A4") "2025-12-13"
Set-TextValue $ws.Range("E4") "0.07"

# Row 5 - Black Lotus
Set-TextValue $ws.Range("A5") "2025-12-13"

# Row 6 - Mox Pearl
Set-TextValue $ws.Range("A6") "2025-12-13"

# Row 7 - Totals
Set-TextValue $ws.Range("E7") "2.68"
Set-TextValue $ws.Range("F7") "Total P/L: `$-5.39"
